$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-121 down to 55-122
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with its data
$ws.Cells.Item(54, 1).Value = 4
$ws.Cells.Item(54, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(54, 3).Value = "Los Lagos"
$ws.Cells.Item(54, 4).Value = 44483
$ws.Cells.Item(54, 5).Value = 10
$ws.Cells.Item(54, 6).Value = 100112039
$ws.Cells.Item(54, 7).Value = "Ciboulette"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 120
$ws.Cells.Item(54, 11).Value = 2500
$ws.Cells.Item(54, 12).Value = 2500
$ws.Cells.Item(54, 13).Value = 2500
$ws.Cells.Item(54, 14).Value = "`$/docena de atados"
$ws.Cells.Item(54, 15).Value = "Región Metropolitana"
$ws.Cells.Item(54, 16).Value = 833
$ws.Cells.Item(54, 17).Value = 3
$ws.Cells.Item(54, 18).Value = "Hortaliza"
